$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The first four data rows of the original sheet (rows 2-5) are removed;
# everything below shifts up so the old row 6 becomes the new row 2, and
# so on, shrinking the table from 13 data rows down to 9. Deleting whole
# rows 2:5 both removes that data and shifts the remaining rows upward,
# matching the target sheet exactly.
$ws.Rows("2:5").Delete()

# Refresh the selection to match the post-edit state.
$ws.Range("A2:XFD4").Select() | Out-Null

# Re-apply the sort over the shrunk range so the sheet's remembered sort
# state matches the new extent instead of pointing past the end of the
# data (it previously covered A2:H13 / A5:A13).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A4:A9")) | Out-Null
$ws.Sort.SetRange($ws.Range("A2:H9"))
$ws.Sort.Header = 0
$ws.Sort.Apply()
